# Fixed Los Angeles Clippers & Miami Heat Logic
# Updates the Elo ratings table + game-prediction section with recalculated
# values, swaps the Clippers/Mavericks rank order, and widens the win
# probability number format so it isn't rounded to 3 decimals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Elo rating table (rows 5-34), column C -------------------------------
$ws.Range("C5").Value  = 1797.1
$ws.Range("C6").Value  = 1777.48
$ws.Range("C7").Value  = 1683.9
$ws.Range("C8").Value  = 1661.8
$ws.Range("C9").Value  = 1616.82
$ws.Range("C10").Value = 1612.74
$ws.Range("C11").Value = 1601.3
$ws.Range("C12").Value = 1589.98
$ws.Range("C13").Value = 1588.03
$ws.Range("C14").Value = 1582.82
$ws.Range("C15").Value = 1578.14
$ws.Range("C16").Value = 1551.64
$ws.Range("C17").Value = 1541.98
$ws.Range("C18").Value = 1529.16
$ws.Range("C19").Value = 1502.88

# Row 20 / 21 swap: Dallas Mavericks now ranks above Los Angeles Clippers.
$ws.Range("B20").Value = "Dallas Mavericks"
$ws.Range("C20").Value = 1478.74
$ws.Range("B21").Value = "Los Angeles Clippers"
$ws.Range("C21").Value = 1470.47

$ws.Range("C22").Value = 1465.62
$ws.Range("C23").Value = 1438.25
$ws.Range("C24").Value = 1438.01
$ws.Range("C25").Value = 1436.24
$ws.Range("C27").Value = 1388.17
$ws.Range("C28").Value = 1378.24
$ws.Range("C29").Value = 1358.1
$ws.Range("C30").Value = 1355.48
$ws.Range("C31").Value = 1325.78
$ws.Range("C32").Value = 1285.63
$ws.Range("C33").Value = 1277.19
$ws.Range("C34").Value = 1262.75

# --- Game predictions (rows 39-42) -----------------------------------------
$ws.Range("C39").Value = 1797.1
$ws.Range("D39").Value = 1355.48
$ws.Range("F39").Value = 0.9576216581651262

$ws.Range("C40").Value = 1578.14
$ws.Range("D40").Value = 1277.19
$ws.Range("F40").Value = 0.9095418528145847

$ws.Range("C41").Value = 1551.64
$ws.Range("D41").Value = 1541.98
$ws.Range("F41").Value = 0.6527736812238141

$ws.Range("C42").Value = 1325.78
$ws.Range("D42").Value = 1470.47
$ws.Range("F42").Value = 0.4360382848393158

# --- Number format: show more decimal precision on win-probability column --
$ws.Range("F39:F42").NumberFormat = "0.############"
